# Insert a new weekly data row for "Haba" (Feria Lagunitas de Puerto Montt)
# at sheet row 63, pushing the existing rows 63-121 down to 64-122.
#
# New row 63 data:
#   Fecha (D) = 44893, Precio min/max/prom (K/L/M) = 16000,
#   Origen (O) = "Región del Maule", Precio $/Kg (P) = 640
# All other columns repeat the constant per-market values used throughout
# the sheet (A, B, C, E, F, G, H, I, J, N, Q, R).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 63 (and everything below it) down by one row, copying the
# formatting of the row above (this is what Excel's own Insert does, and
# it is what carries the date-format style (s="2") onto the new D63 cell).
$ws.Rows.Item(63).Insert()

$ws.Range("A63").Value = 4
$ws.Range("B63").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C63").Value = "Los Lagos"
$ws.Range("D63").Value = 44893
$ws.Range("E63").Value = 10
$ws.Range("F63").Value = 100112026
$ws.Range("G63").Value = "Haba"
$ws.Range("H63").Value = "Sin especificar"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 80
$ws.Range("K63").Value = 16000
$ws.Range("L63").Value = 16000
$ws.Range("M63").Value = 16000
$ws.Range("N63").Value = "`$/saco 25 kilos"
$ws.Range("O63").Value = "Región del Maule"
$ws.Range("P63").Value = 640
$ws.Range("Q63").Value = 25
$ws.Range("R63").Value = "Hortaliza"
